$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10826
$ws1.Range("F4").Value = 68
$ws1.Range("F5").Value = 734

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 7

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10826
$ws4.Range("F4").Value = 68
$ws4.Range("F5").Value = 734
$ws4.Range("F6").Value = 7
